$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-294) holds the "Förändrad" (Changed) date.
# Update the date serial value from 45171 (2023-09-02) to 45172 (2023-09-03)
# for every data row.
$ws.Range("C2:C294").Value = 45172
